$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row with the food cost item
$ws.Range("A7").Value = "Comida"
$ws.Range("B7").Value = 1000

# Match the currency style used by the other value cells (B5, B6)
$ws.Range("B7").NumberFormat = $ws.Range("B6").NumberFormat

# Update selection to match the post-edit active cell
$ws.Range("E7").Select()

$wb.Save()
